# chore: adapt column header formatting to respective input file names
#
# The "AHB-Diff" sheet has 21 columns: 10 "old" AHB columns (A:J), a
# "diff" column (K), and 10 "new" AHB columns (L:U). The column headers
# used the generic suffixes "_old"/"_new"; rename them to carry the
# concrete format-version suffixes ("_FV2310" / "_FV2404") instead, wrap
# the sheet's used range in a native Excel table ("Table1") and freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base header names (without the old "_old"/"_new" suffix), in column order.
$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1-10): "<name>_old" -> "<name>_FV2310"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseHeaders[$i])_FV2310"
}

# Column K (11) is "diff" and is left untouched.

# Columns L..U (12-21): "<name>_new" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseHeaders[$i])_FV2404"
}

# Turn the used range A1:U66 into a native table named "Table1" with a
# header row (matches the new <tableParts>/xl/tables/table1.xml part).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split under row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Renamed headers, added Table1 over A1:U66, froze top row."
